$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.324.90"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "1.651.66"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.09"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.511"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.44"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0894"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").Value = "1.885.57"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "1.647.10"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("E15").Value = "  +3.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.43"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "27.333.27"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.59"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -6.21%  "
$ws.Range("D19").Value = "0.0₃0727"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.40"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.03"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("E31").Value = "  -4.26%  "
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").Value = "1.428.73"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.13"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.906"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.571"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.03"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.21"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").Value = "1.794.23"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.95"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.68%  "
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.75"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.59%  "
